# SO_Test.xlsx update
# - Row 9: mark COMPLETED and RECEIVED as YES (with date/clerk), bump # OF CALLS
#          and record CALL1 info ("picked up order was placed" label)
# - Rows 12, 16, 21: mark COMPLETED as YES with date/clerk (order was picked up)
# - Row 28: order fully completed -> results frame refreshed / row cleared

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$today = "03/14/2024"
$clerk = "abake"

# Helper: write a date-like string as plain text, matching the other date
# columns in this sheet which already store dates as text (e.g. "03/13/2024")
# rather than as real Excel date serials.
function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# --- Row 9 ---
$ws.Range("U9").Value  = "YES"
Set-TextValue "V9" $today
$ws.Range("W9").Value  = $clerk
$ws.Range("X9").Value  = "YES"
Set-TextValue "Y9" $today
$ws.Range("Z9").Value  = $clerk
$ws.Range("AA9").Value = 1
Set-TextValue "AB9" $today
$ws.Range("AC9").Value = $clerk
$ws.Range("AD9").Value = "picked up"

# --- Row 12 ---
$ws.Range("U12").Value = "YES"
Set-TextValue "V12" $today
$ws.Range("W12").Value = $clerk

# --- Row 16 ---
$ws.Range("U16").Value = "YES"
Set-TextValue "V16" $today
$ws.Range("W16").Value = $clerk

# --- Row 21 ---
$ws.Range("U21").Value = "YES"
Set-TextValue "V21" $today
$ws.Range("W21").Value = $clerk

# --- Row 28: order completed, clear/refresh the row ---
$row28Cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ")
foreach ($col in $row28Cols) {
    $ws.Range($col + "28").Value = " "
}
